# Commit: "added gecko driver; added string split to capture 7th page in results"
#
# - startup!G2 (location_path): point at the new exported data file
#   (was a stale Trial_Test1Data.xlsx path, now the fresh Neo4j export).
# - testcase!A3:A6 (Browser): switch the browser used to run the test
#   cases from Chrome to Firefox (i.e. "added gecko driver" -> geckodriver
#   is the WebDriver used for Firefox).

$wb = $excel.ActiveWorkbook

$startup  = $wb.Worksheets.Item("startup")
$testcase = $wb.Worksheets.Item("testcase")

# startup sheet: update the location_path value
$startup.Range("G2").Value = "C:\Users\radhakrishnang2\Desktop\DataCommons_Automation\CTDC_Automation\TestData\DatafromNeo4j.xlsx"

# testcase sheet: Chrome -> Firefox for every test-case row
$testcase.Range("A3").Value = "Firefox"
$testcase.Range("A4").Value = "Firefox"
$testcase.Range("A5").Value = "Firefox"
$testcase.Range("A6").Value = "Firefox"

# Leave the selection/active-cell state the way the author's save shows it
$startup.Activate()
$startup.Range("G2").Select()

$testcase.Activate()
$testcase.Range("A1").Select()
$testcase.Range("A12").Select()
